$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    (Heading1 title) paragraph.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaRange.Text = "Meta description: Experience the ancient Egypt theme with A While On The Nile and enjoy exciting bonus features - play for free and potentially win big."

# Re-fetch the paragraph's range and bold just the "Meta description" lead-in.
$metaFull = $metaPara.Range
$boldRange = $d.Range($metaFull.Start, $metaFull.Start + 16)
$boldRange.Font.Bold = 1

# ------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the
#    document ("Play A While On The Nile Slot for Free - Exciting
#    Bonus Features"). Search from the end (skipping the real title
#    in paragraph 1) since Paragraph.Range.Text includes the trailing
#    paragraph mark and must be trimmed before comparing.
# ------------------------------------------------------------------
$dupParaIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $pp = $d.Paragraphs($i)
    if ($pp.Range.Text.Trim() -eq "Play A While On The Nile Slot for Free - Exciting Bonus Features") {
        $dupParaIndex = $i
        break
    }
}
if ($dupParaIndex -ge 2) {
    $d.Paragraphs($dupParaIndex).Range.Delete()
}

# ------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    feature-image prompt copy, keeping its run formatting intact.
#    (Direct Range.Text assignment is used instead of Find/Replace so
#    that straight quotes survive instead of being "smart quoted".)
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastTrimmed = $d.Range($lastRange.Start, $lastRange.End - 1)
$lastTrimmed.Text = "Create a Feature Image Prompt: Design a cartoon-style feature image for ""A While On The Nile"" online slot game that showcases a happy Maya warrior wearing glasses. The warrior should be depicted in an Egyptian-themed outfit, possibly holding a tablet with hieroglyphics or standing by the Nile river. The image should have bright colors and should be eye-catching to potential players. Make sure to include the game's title in the image prominently."

Write-Output "done"
